$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.784.43'
$ws.Range('E2').Value = '  +1.76%  '
$ws.Range('D3').Value = '1.879.85'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '332.85'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +2.81%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.003'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.13%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4723'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +4.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3956'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.34%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '47.89'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08062'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +1.72%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.033'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.77%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.29'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.23%  '
$ws.Range('D13').Value = '1.879.46'
$ws.Range('E13').Value = '  +1.40%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.987'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.35%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.133'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -0.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.006'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +0.40%  '
$ws.Range('B17').Value = 'Litecoin'
$ws.Range('C17').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '87.28'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +1.63%  '
$ws.Range('B18').Value = 'ShibaInu'
$ws.Range('C18').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001046'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06665'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +1.71%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.30'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.42%  '
$ws.Range('E21').Value = '  +0.10%  '
$ws.Range('D22').Value = '27.797.75'
$ws.Range('E23').Value = '  +0.63%  '
$ws.Range('E24').Value = '  +1.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.307'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.91%  '
$ws.Range('D26').Value = '2.080.86'
$ws.Range('E26').Value = '  +0.51%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '159.49'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.76%  '
$ws.Range('E28').Value = '  +1.57%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.108'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +2.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '5.584'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +2.10%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '121.96'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.9829'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +5.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09543'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.51%  '
$ws.Range('E34').Value = '  -0.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.595'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +0.32%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.366'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.80%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.06132'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.21%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.02259'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.52%  '
$ws.Range('E39').Value = '  +1.02%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.137'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.80%  '
$ws.Range('E41').Value = '  +2.03%  '
$ws.Range('E42').Value = '  +0.13%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1905'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.08%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.28'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +1.39%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.275'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.77%  '
$ws.Range('E46').Value = '  +2.11%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.23'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +1.87%  '
$ws.Range('E48').Value = '  +1.76%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.373'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06916'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.60%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '114.42'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +5.24%  '
